$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-18 21:18:45"
$ws.Range("E3").Value = "2026-02-18 21:18:47"
$ws.Range("H3").Value = "71%"
$ws.Range("E4").Value = "2026-02-18 21:18:50"
$ws.Range("H4").Value = "72%"
$ws.Range("J4").Value = "1012.8 hPa"
$ws.Range("O4").Value = "12.2 °C"
$ws.Range("E5").Value = "2026-02-18 21:18:53"
$ws.Range("I5").Value = "1.1 mm"
$ws.Range("E6").Value = "2026-02-18 21:18:55"
$ws.Range("H6").Value = "73%"
$ws.Range("J6").Value = "1012.5 hPa"
$ws.Range("O6").Value = "12.2 °C"
$ws.Range("E7").Value = "2026-02-18 21:18:58"
$ws.Range("H7").Value = "74%"
$ws.Range("J7").Value = "1014.0 hPa"
$ws.Range("O7").Value = "13.8 °C"
$ws.Range("E8").Value = "2026-02-18 21:19:01"
$ws.Range("J8").Value = "1013.8 hPa"
$ws.Range("O8").Value = "11.2 °C"
$ws.Range("E9").Value = "2026-02-18 21:19:03"
$ws.Range("O9").Value = "11.1 °C"
$ws.Range("E10").Value = "2026-02-18 21:19:06"
$ws.Range("E11").Value = "2026-02-18 21:19:09"
$ws.Range("E12").Value = "2026-02-18 21:19:12"
$ws.Range("H12").Value = "86%"
$ws.Range("E13").Value = "2026-02-18 21:19:13"
$ws.Range("H13").Value = "73%"
$ws.Range("J13").Value = "1015.2 hPa"
$ws.Range("E14").Value = "2026-02-18 21:19:14"
$ws.Range("E15").Value = "2026-02-18 21:19:15"
$ws.Range("O15").Value = "10.5 °C"
$ws.Range("E16").Value = "2026-02-18 21:19:16"
$ws.Range("G16").Value = "71 cm"
$ws.Range("H16").Value = "53%"
$ws.Range("I16").Value = "1.0 mm"
$ws.Range("O16").Value = "-0.1 °C"
$ws.Range("E17").Value = "2026-02-18 21:19:17"
$ws.Range("E18").Value = "2026-02-18 21:19:19"
$ws.Range("J18").Value = "1013.0 hPa"
$ws.Range("O18").Value = "11.8 °C"
$ws.Range("E19").Value = "2026-02-18 21:19:20"
$ws.Range("E20").Value = "2026-02-18 21:19:21"
$ws.Range("E21").Value = "2026-02-18 21:19:22"
$ws.Range("J21").Value = "1014.6 hPa"
$ws.Range("O21").Value = "6.6 °C"
$ws.Range("E22").Value = "2026-02-18 21:19:25"
$ws.Range("H22").Value = "91%"
$ws.Range("I22").Value = "1.2 mm"
$ws.Range("L22").Value = "65.5 km/h - 309º 20:59 TU"
$ws.Range("O22").Value = "-1.7 °C"
$ws.Range("E23").Value = "2026-02-18 21:19:28"
$ws.Range("H23").Value = "57%"
$ws.Range("O23").Value = "0.0 °C"
$ws.Range("E24").Value = "2026-02-18 21:19:30"
$ws.Range("J24").Value = "1014.6 hPa"
$ws.Range("L24").Value = "24.8 km/h - 263º 20:58 TU"
$ws.Range("E25").Value = "2026-02-18 21:19:33"
$ws.Range("H25").Value = "48%"
$ws.Range("E26").Value = "2026-02-18 21:19:36"
$ws.Range("H26").Value = "70%"
$ws.Range("J26").Value = "1011.9 hPa"
$ws.Range("E27").Value = "2026-02-18 21:19:38"
$ws.Range("H27").Value = "57%"
$ws.Range("O27").Value = "1.5 °C"
$ws.Range("E28").Value = "2026-02-18 21:19:41"
$ws.Range("J28").Value = "1012.6 hPa"
$ws.Range("E29").Value = "2026-02-18 21:19:44"
$ws.Range("E30").Value = "2026-02-18 21:19:47"
$ws.Range("J30").Value = "1012.2 hPa"
$ws.Range("E31").Value = "2026-02-18 21:19:50"
$ws.Range("H31").Value = "73%"
$ws.Range("J31").Value = "1011.0 hPa"
$ws.Range("O31").Value = "12.6 °C"
$ws.Range("E32").Value = "2026-02-18 21:19:52"
$ws.Range("O32").Value = "8.0 °C"
$ws.Range("E33").Value = "2026-02-18 21:19:55"
$ws.Range("J33").Value = "1013.8 hPa"
$ws.Range("E34").Value = "2026-02-18 21:19:57"
$ws.Range("E35").Value = "2026-02-18 21:20:00"
$ws.Range("J35").Value = "1014.2 hPa"
$ws.Range("E36").Value = "2026-02-18 21:20:03"
$ws.Range("H36").Value = "85%"
$ws.Range("J36").Value = "1012.7 hPa"
$ws.Range("E37").Value = "2026-02-18 21:20:06"
$ws.Range("J37").Value = "1014.3 hPa"
$ws.Range("O37").Value = "6.2 °C"
$ws.Range("E38").Value = "2026-02-18 21:20:08"
$ws.Range("E39").Value = "2026-02-18 21:20:11"
$ws.Range("E40").Value = "2026-02-18 21:20:14"
$ws.Range("I40").Value = "0.5 mm"
$ws.Range("J40").Value = "1015.2 hPa"
$ws.Range("E41").Value = "2026-02-18 21:20:16"
$ws.Range("J41").Value = "1014.3 hPa"
$ws.Range("E42").Value = "2026-02-18 21:20:19"
$ws.Range("H42").Value = "84%"
$ws.Range("E43").Value = "2026-02-18 21:20:21"
$ws.Range("H43").Value = "78%"
$ws.Range("E44").Value = "2026-02-18 21:20:23"
$ws.Range("H44").Value = "74%"
$ws.Range("E45").Value = "2026-02-18 21:20:26"
$ws.Range("H45").Value = "64%"
$ws.Range("I45").Value = "0.7 mm"
$ws.Range("J45").Value = "1011.6 hPa"
$ws.Range("O45").Value = "7.4 °C"
$ws.Range("E46").Value = "2026-02-18 21:20:29"
$ws.Range("H46").Value = "84%"
$ws.Range("J46").Value = "1014.7 hPa"
